$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(16, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1607', 24640, 781242)
    ,@(17, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1608', 24640, 781242)
    ,@(18, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1609', 24640, 781242)
    ,@(19, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1610', 24640, 781242)
    ,@(20, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1611', 24640, 781242)
    ,@(21, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1612', 24640, 781242)
    ,@(22, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1701', 24640, 781242)
    ,@(23, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1702', 24640, 781242)
    ,@(24, '1101442161', 'MARIA BERNARDA GUTIERREZ PEREZ', '1703', 29509, 737717)
    ,@(25, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1703', 24640, 781242)
    ,@(26, '45530618', 'SANDRA PATRICIA CONTRERAS ROMERO', '1703', 29509, 781242)
    ,@(27, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1704', 24640, 781242)
    ,@(28, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1705', 24640, 781242)
    ,@(29, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1706', 24640, 781242)
    ,@(30, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1707', 24640, 781242)
    ,@(31, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1708', 24640, 781242)
    ,@(32, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1709', 24640, 781242)
    ,@(33, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1710', 24640, 781242)
    ,@(34, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1711', 24640, 781242)
    ,@(35, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1712', 24640, 781242)
    ,@(36, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1801', 24640, 781242)
    ,@(37, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1802', 24640, 781242)
    ,@(38, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1803', 24640, 781242)
    ,@(39, '1143357077', 'CINDY STEPHANIE PETRO LARA', '1804', 31249, 781242)
    ,@(40, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1804', 24640, 781242)
    ,@(41, '1143357077', 'CINDY STEPHANIE PETRO LARA', '1805', 31249, 781242)
    ,@(42, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1805', 24640, 781242)
    ,@(43, '1143357077', 'CINDY STEPHANIE PETRO LARA', '1806', 31249, 781242)
    ,@(44, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1806', 24640, 781242)
    ,@(45, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1807', 24640, 781242)
    ,@(46, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1808', 24640, 781242)
    ,@(47, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1809', 31249, 781242)
    ,@(48, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1810', 31249, 781242)
    ,@(49, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1811', 31249, 781242)
    ,@(50, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1812', 31249, 781242)
    ,@(51, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1901', 31249, 781242)
    ,@(52, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1902', 31249, 781242)
    ,@(53, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1903', 31249, 781242)
    ,@(54, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1904', 31249, 781242)
    ,@(55, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1905', 31249, 781242)
    ,@(56, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1906', 31249, 781242)
    ,@(57, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1907', 31249, 781242)
    ,@(58, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1908', 31249, 781242)
    ,@(59, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1909', 31249, 781242)
    ,@(60, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1910', 31249, 781242)
    ,@(61, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1911', 31249, 781242)
    ,@(62, '22725529', 'DELIA ROSA MANOTAS TOVAR', '1912', 31249, 781242)
    ,@(63, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2001', 31249, 781242)
    ,@(64, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2002', 31249, 781242)
    ,@(65, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2003', 31249, 781242)
    ,@(66, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2004', 31249, 781242)
    ,@(67, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2005', 31249, 781242)
    ,@(68, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2006', 31249, 781242)
    ,@(69, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2007', 31249, 781242)
    ,@(70, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2008', 31249, 781242)
    ,@(71, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2009', 31249, 781242)
    ,@(72, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2010', 31249, 781242)
    ,@(73, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2011', 31249, 781242)
    ,@(74, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2012', 31249, 781242)
    ,@(75, '45452537', 'NELVA ROSA RIVERA ZABALETA', '2101', 105334, 3619371)
    ,@(76, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2101', 31249, 781242)
    ,@(77, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2102', 31249, 781242)
    ,@(78, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2103', 31249, 781242)
    ,@(79, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2104', 31249, 781242)
    ,@(80, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2105', 31249, 781242)
    ,@(81, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2106', 31249, 781242)
    ,@(82, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2107', 31249, 781242)
    ,@(83, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2108', 31249, 781242)
    ,@(84, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2109', 31249, 781242)
    ,@(85, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2110', 31249, 781242)
    ,@(86, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2111', 31249, 781242)
    ,@(87, '22725529', 'DELIA ROSA MANOTAS TOVAR', '2112', 8333, 781242)
    ,@(88, '1049827118', 'MARY ANTONIA MOYA RUBIO', '2403', 52000, 1423500)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = [string]$row[1]
    $ws.Cells.Item($r, 4).Value = [string]$row[2]
    $ws.Cells.Item($r, 5).Value = [string]$row[3]
    $ws.Cells.Item($r, 6).Value = [double]$row[4]
    $ws.Cells.Item($r, 7).Value = [double]$row[5]
}
